# Adds two new columns, I ("I0") and J ("IF"), to the worksheet, with
# per-row numeric values, matching the "I0 and IF added" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same bold / centered / bordered style
# (style index "1" in the original workbook) before setting their text.
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(1, 10).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data rows (2-38) -------------------------------------------------
# I column ("I0") values per row (rows 2..38)
$iVals = @(1,6,1,7,1,1,1,1,1,1,1,1,1,2,1,1,1,1,1,1,2,1,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
# J column ("IF") values per row (rows 2..38)
$jVals = @(6,7,4,9,4,4,5,7,8,6,7,6,5,7,6,7,6,6,6,7,6,6,7,5,6,5,6,6,6,6,6,5,6,5,4,3,1)

for ($row = 2; $row -le 38; $row++) {
    $idx = $row - 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
